# Reorder character-style rPr children so <w:b/>/<w:i/> come before
# <w:color/>, matching wml.xsd's CT_RPr sequence (fixes the
# OOXMLValidator "Sch_UnexpectedElementContentExpectingComplex" warning
# on KeywordTok and friends). Touching Font.Bold / Font.Italic on the
# style causes the rPr to be rebuilt in schema order.

$d = $word.ActiveDocument

$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnly) {
    $style = $d.Styles($styleName)
    $style.Font.Bold = $true
}

$italicOnly = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnly) {
    $style = $d.Styles($styleName)
    $style.Font.Italic = $true
}

$boldAndItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldAndItalic) {
    $style = $d.Styles($styleName)
    $style.Font.Bold = $true
    $style.Font.Italic = $true
}
